# Update the Ghrh-Vipr1 NATMI LR-pairs sheet with refreshed TPM-derived
# specificity numbers, keep only the "MuSCs" sending-cluster rows (the old
# "ECs" sending-cluster rows are dropped), and drop the now-unused rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Sending cluster MuSCs -> Target cluster ECs -----------------
$ws.Range("A2").Value2 = "MuSCs"
$ws.Range("B2").Value2 = "Ghrh"
$ws.Range("C2").Value2 = "Vipr1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.09527133333333333
$ws.Range("H2").Value2 = 0.285814
$ws.Range("I2").Value2 = 1
$ws.Range("J2").Value2 = 1
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.8059226666666667
$ws.Range("N2").Value2 = 2.417768
$ws.Range("O2").Value2 = 0.1314814101815314
$ws.Range("P2").Value2 = 0.1314814101815314
$ws.Range("Q2").Value2 = 0.07678132701688889
$ws.Range("R2").Value2 = 0.6910319431520001
$ws.Range("S2").Value2 = 0.1314814101815314
$ws.Range("T2").Value2 = 0.1314814101815314

# ---- Row 3: Sending cluster MuSCs -> Target cluster MuSCs ---------------
$ws.Range("A3").Value2 = "MuSCs"
$ws.Range("B3").Value2 = "Ghrh"
$ws.Range("C3").Value2 = "Vipr1"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.09527133333333333
$ws.Range("H3").Value2 = 0.285814
$ws.Range("I3").Value2 = 1
$ws.Range("J3").Value2 = 1
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 5.323633333333333
$ws.Range("N3").Value2 = 15.9709
$ws.Range("O3").Value2 = 0.8685185898184686
$ws.Range("P3").Value2 = 0.8685185898184687
$ws.Range("Q3").Value2 = 0.5071896458444445
$ws.Range("R3").Value2 = 4.5647068126
$ws.Range("S3").Value2 = 0.8685185898184686
$ws.Range("T3").Value2 = 0.8685185898184687

# ---- Drop the two now-obsolete rows (old "MuSCs" duplicate block) -------
# Delete from the bottom up so row numbers of earlier rows stay stable.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
